$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-18
$values = @(
    @(6, 8),
    @(6, 9),
    @(8, 9),
    @(4, 4),
    @(6, 7),
    @(5, 6),
    @(5, 7),
    @(3, 6),
    @(7, 7),
    @(5, 7),
    @(6, 7),
    @(8, 9),
    @(8, 9),
    @(8, 8),
    @(7, 7),
    @(4, 5),
    @(5, 6)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
    $row++
}
